# Weekly driver report update for 2025-05-05
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Bad Drivers" summary table (row 3 = data row, row 4 = totals row)
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 98.5
$ws.Range("C4").Value = 13

# "Good Drivers" detail table (row 12)
$ws.Range("B12").Value = 1074341
$ws.Range("C12").Value = 4452
$ws.Range("E12").Value = 1647
$ws.Range("F12").Value = 1078973
